# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, matching the refreshed export.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 218
$ws1.Range("F5").Value = 1345
$ws1.Range("F7").Value = 598
$ws1.Range("F8").Value = 121
$ws1.Range("F9").Value = 580
$ws1.Range("F10").Value = 25
$ws1.Range("F11").Value = 655
$ws1.Range("F14").Value = 153
$ws1.Range("F15").Value = 225

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 13
$ws2.Range("F6").Value = 3

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6286
$ws3.Range("F4").Value = 1920

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6286
$ws4.Range("F4").Value = 1920
$ws4.Range("F9").Value = 13
$ws4.Range("F11").Value = 218
$ws4.Range("F12").Value = 3
$ws4.Range("F15").Value = 1345
$ws4.Range("F20").Value = 598
$ws4.Range("F22").Value = 121
$ws4.Range("F23").Value = 580
$ws4.Range("F24").Value = 25
$ws4.Range("F26").Value = 655
$ws4.Range("F31").Value = 153
$ws4.Range("F37").Value = 225
